$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(24,34,91,25,35,31),
    @(41,43,98,22,83,52),
    @(15,27,94,17,39,30),
    @(15,19,91,25,58,55),
    @(41,66,93,47,68,75),
    @(10,73,37,11,77,13),
    @(9,81,82,73,88,68),
    @(33,51,52,4,77,43),
    @(27,81,90,65,95,72),
    @(22,45,84,78,86,50),
    @(51,48,84,4,93,76),
    @(17,32,86,15,93,45),
    @(12,35,46,23,71,36),
    @(2,20,50,19,90,36),
    @(8,3,62,54,92,57),
    @(15,32,79,47,80,51),
    @(11,69,36,33,96,88),
    @(44,45,66,16,97,81),
    @(32,60,77,52,91,81),
    @(5,12,42,39,84,67),
    @(4,13,19,7,84,99)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $excelRow = $startRow + $i
    for ($j = 0; $j -lt $row.Length; $j++) {
        $col = $j + 2  # column B is index 2
        $ws.Cells.Item($excelRow, $col).Value = $row[$j]
    }
}

$wb.Save()
